# Apply updated bus voltage magnitude (vm_pu) results for the 380 kV case
# (B column = slack/reference bus voltage changed from 1.05 to 1.02 p.u.,
#  all dependent bus voltages recomputed accordingly).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number, then values for columns B,C,D,E,F,I,J,K,L,M,N (column H has no data)
$rowsData = @(
    ,@(2, 1.02, 1.042578163836886, 1.049719097510942, 1.04623855066462, 1.057125030866636, 1.041890428889587, 1.047653232839698, 1.052475373860962, 1.049004545774487, 1.059860858288384, 1.049141020597836)
    ,@(3, 1.02, 1.043543760084638, 1.05050145029989, 1.047153034049339, 1.058088752324795, 1.042141689228796, 1.048265244721133, 1.053070106377957, 1.049730363953588, 1.060637973089561, 1.049753901606341)
    ,@(4, 1.02, 1.044169020690681, 1.051008079850108, 1.047745553942844, 1.058713191254532, 1.042303308729274, 1.04866109119299, 1.053454673883419, 1.050200173777044, 1.061141041811192, 1.050150310225592)
    ,@(5, 1.02, 1.044431988438199, 1.051221160274909, 1.047994836798155, 1.058975906891709, 1.042371022611434, 1.04882746444881, 1.053616281930046, 1.050397718439811, 1.061352584372667, 1.050316919750519)
    ,@(6, 1.02, 1.044476148201146, 1.051256942855838, 1.048036703447562, 1.059020029785336, 1.042382378513123, 1.048855396845699, 1.053643412845816, 1.050430889177203, 1.061388106346583, 1.050344891814616)
    ,@(7, 1.02, 1.044172534054477, 1.051010926675823, 1.047748884137283, 1.058716700883419, 1.042304214433566, 1.048663314441141, 1.053456833551145, 1.050202813235824, 1.061143868248363, 1.05015253663101)
    ,@(8, 1.02, 1.042904396859624, 1.049983414795209, 1.046547440751524, 1.057450549012656, 1.041975542553763, 1.047860099023013, 1.052676420673303, 1.049249806358226, 1.060123440980103, 1.049348180554859)
    ,@(9, 1.02, 1.040673299344784, 1.048175885520555, 1.044436430100531, 1.055225965012117, 1.041389028231715, 1.046443494901375, 1.051299248956834, 1.047571730550132, 1.058327074880324, 1.047929564692882)
    ,@(10, 1.02, 1.039188315477112, 1.046973005900648, 1.043033249384586, 1.053747374271417, 1.040993103125866, 1.045498309554476, 1.05037985194898, 1.046453905771584, 1.057130739318784, 1.046983037074376)
    ,@(11, 1.02, 1.038545883024345, 1.04645266790083, 1.04242665724037, 1.053108201658177, 1.040820502990185, 1.045088858161867, 1.049981449347026, 1.045970098479334, 1.056613019800207, 1.046573004213824)
    ,@(12, 1.02, 1.038307342133198, 1.046259469990352, 1.042201492390835, 1.052870945907459, 1.040756217471493, 1.044936743414892, 1.049833420985721, 1.045790424581191, 1.05642076200481, 1.046420673446461)
    ,@(13, 1.02, 1.038358506010403, 1.04630090801413, 1.04224978418829, 1.052921830748865, 1.040770014806942, 1.044969373720996, 1.049865175544923, 1.045828963704266, 1.056461999844575, 1.046453350091342)
    ,@(14, 1.02, 1.038526163378541, 1.046436696487346, 1.042408041967583, 1.053088586724451, 1.040815192679146, 1.045076284840908, 1.04996921416722, 1.045955245886828, 1.056597126761083, 1.046560413037307)
    ,@(15, 1.02, 1.038629474182085, 1.046520370722094, 1.042505569781438, 1.053191352007591, 1.040843005194416, 1.045142152852261, 1.05003331000403, 1.04603305699577, 1.056680389076174, 1.046626374588792)
    ,@(16, 1.02, 1.0392309637837, 1.047007550006406, 1.043073527950483, 1.053789816627493, 1.041004533583802, 1.045525479769347, 1.050406286416733, 1.0464860191414, 1.057165105096152, 1.047010245874069)
    ,@(17, 1.02, 1.039608416841679, 1.047313283874171, 1.04343006009266, 1.054165504112274, 1.041105545287832, 1.045765882987739, 1.050640165560218, 1.046770209375996, 1.057469236094537, 1.047250990492604)
    ,@(18, 1.02, 1.039828634292201, 1.047491663054766, 1.043638115127814, 1.054384739134333, 1.041164351565492, 1.04590608870893, 1.050776554535578, 1.04693599372648, 1.057646659602914, 1.047391395322002)
    ,@(19, 1.02, 1.039903732139458, 1.047552494155717, 1.043709072746063, 1.054459510025469, 1.041184383952594, 1.045953892260727, 1.050823054746419, 1.046992525439851, 1.057707161331118, 1.047439266760325)
    ,@(20, 1.02, 1.039567913958452, 1.047280476365481, 1.043391797643948, 1.05412518575234, 1.04109471928406, 1.045740091792724, 1.050615075511624, 1.046739716275348, 1.057436602705698, 1.047225162671135)
    ,@(21, 1.02, 1.038476790051738, 1.046396707982209, 1.042361434833271, 1.053039476780996, 1.040801893725574, 1.045044802899824, 1.049938578591199, 1.045918057994917, 1.056557333966679, 1.046528886388256)
    ,@(22, 1.02, 1.037791261573243, 1.045841504149601, 1.041714476018978, 1.052357782925859, 1.040616775348902, 1.044607495570707, 1.049512984565702, 1.045401644078075, 1.056004771083138, 1.046090958032572)
    ,@(23, 1.02, 1.038154625146166, 1.046135784559336, 1.042057358185147, 1.052719072784366, 1.040715005442472, 1.044839334576365, 1.049738623771523, 1.045675386061198, 1.056297669371526, 1.04632312627621)
    ,@(24, 1.02, 1.039586215296582, 1.047295300512744, 1.043409086501404, 1.054143403565955, 1.041099611435687, 1.045751745777423, 1.050626412713557, 1.046753494730165, 1.057451348239763, 1.047236833205829)
    ,@(25, 1.02, 1.041249669374877, 1.048642803155721, 1.044981448746088, 1.055800291395286, 1.041541524691024, 1.046809862511313, 1.051655510707519, 1.048005400181039, 1.058791264128797, 1.048296452586849)
)

$columns = @("B","C","D","E","F","I","J","K","L","M","N")

foreach ($entry in $rowsData) {
    $rowNum = $entry[0]
    for ($i = 0; $i -lt $columns.Count; $i++) {
        $colLetter = $columns[$i]
        $value = $entry[$i + 1]
        $ws.Range("$colLetter$rowNum").Value = $value
    }
}

Write-Output "Updated $($rowsData.Count) rows x $($columns.Count) columns of vm_pu results"
